$wb = $excel.ActiveWorkbook

# --- Text change: "Ready for handoff" -> "In Translation" on every sheet ---
# This status text shows up in the Overview sheet (per-language status
# columns) as well as in each language sheet's "Status" column.
foreach ($ws in $wb.Worksheets) {
    $ws.Cells.Replace("Ready for handoff", "In Translation") | Out-Null
}

# --- Column width change on the columns that held the status text ---
# Overview: columns E (zh-cn) and F (de-de)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E1:F1").ColumnWidth = 12.5

# zh-cn / de-de: column C ("Status")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C1").ColumnWidth = 12.5

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C1").ColumnWidth = 12.5
